$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "66.738.48"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +4.09%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "3.257.13"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +7.18%  "

$ws.Range("E4").Value = "  +0.12%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "584.11"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +5.17%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "154.35"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +8.92%  "

$ws.Range("E7").Value = "  +0.08%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "3.249.36"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +7.21%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.516"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +5.36%  "

$ws.Range("E10").Value = "  +8.74%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.166"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +6.15%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.491"
$cell.Style = "Normal"

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "38.11"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +3.85%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "0.0000236"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +6.19%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "3.783.90"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +7.46%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "560.14"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +13.48%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "66.779.30"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +4.15%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "3.254.12"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +6.94%  "

$ws.Range("E19").Value = "  +3.41%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "7.16"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +6.47%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "14.54"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +5.37%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "0.747"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +8.11%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "7.80"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +8.50%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "13.68"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +7.14%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "82.14"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +3.59%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +0.04%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "9.31"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +18.00%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "3.00"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +8.92%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "2.25"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +6.19%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "27.97"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +6.95%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "2.78"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +5.68%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -0.01%  "

$ws.Range("E33").Value = "  +6.25%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "563.59"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +8.85%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "5.74"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +4.18%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "6.41"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +7.33%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.0461"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +13.72%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "55.45"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +4.95%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.133"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +8.08%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.0867"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +7.41%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "3.07"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +13.99%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "3.178.91"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +9.27%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "8.66"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +2.89%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.277"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +11.96%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "2.34"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +10.52%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "26.54"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +4.62%  "

$ws.Range("E47").Value = "  +0.08%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0560"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +3.84%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "126.32"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +5.24%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.114"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +3.31%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "2.25"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +8.56%  "
